# Update cached "today" text on every datetimeFigureOut date field
# (slide master, all slide layouts, and the notes master) from
# 2023-06-29 to 2023-11-07, and annotate the "OSMTILE" extent example
# on slide 5 with a checked="checked" attribute.

$p = $ppt.ActivePresentation

$oldDate = "2023-06-29"
$newDate = "2023-11-07"

function Update-DatePlaceholders {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }

        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if (-not $isDatePlaceholder) { continue }

        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide master
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout off the (single) slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholders $p.NotesMaster.Shapes

# Slide 5: "<extent units="OSMTILE"> ..." example -- add a
# checked="checked" attribute example right after the closing quote.
$slide5 = $p.Slides.Item(5)
$rect = $slide5.Shapes.Item(40)
if ($rect.Name -ne "Rectangle 84") {
    for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
        if ($slide5.Shapes.Item($i).Name -eq "Rectangle 84") {
            $rect = $slide5.Shapes.Item($i)
            break
        }
    }
}

$tr2 = $rect.TextFrame.TextRange
$fullText = $tr2.Text
$needle = '"OSMTILE"'
$idx = $fullText.IndexOf($needle)
if ($idx -ge 0) {
    $sub = $tr2.Characters($idx + 1, $needle.Length)
    $sub.Text = '"OSMTILE“ checked=“checked”'
}
